$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H9").Value = 323.81818
$ws.Range("I9").Value = 194.71428
$ws.Range("J9").Value = 549.75
$ws.Range("K9").Value = 194.71428
$ws.Range("L9").Value = 549.75
$ws.Range("M9").Value = -25.71428
$ws.Range("N9").Value = -887.75

$ws.Range("H26").Value = 0
$ws.Range("I26").Value = 0
$ws.Range("J26").Value = 0
$ws.Range("K26").Value = 0
$ws.Range("L26").ClearContents()
$ws.Range("N26").Value = 0

$ws.Range("H38").Value = 491.91666
$ws.Range("I38").Value = 491.91666
$ws.Range("J38").Value = 0
$ws.Range("K38").Value = 1475.74998
$ws.Range("L38").Value = 0
$ws.Range("M38").Value = -1103.74998

$ws.Range("H39").Value = 691.8
$ws.Range("I39").Value = 575
$ws.Range("J39").Value = 769.6667
$ws.Range("K39").Value = 1725
$ws.Range("L39").Value = 2309.0001
$ws.Range("M39").Value = -1429
$ws.Range("N39").Value = -2901.0001

$ws.Range("H58").Value = 2572.25
$ws.Range("I58").Value = 96.333336
$ws.Range("J58").Value = 10000
$ws.Range("K58").Value = 289.000008
$ws.Range("L58").Value = 30000
$ws.Range("M58").Value = -139.000008
$ws.Range("N58").Value = -30300

$ws.Range("H62").Value = 1975.4546
$ws.Range("I62").Value = 1966.375
$ws.Range("J62").Value = 1999.6666
$ws.Range("K62").Value = 1966.375
$ws.Range("L62").Value = 1999.6666
$ws.Range("M62").Value = -1342.375
$ws.Range("N62").Value = -3247.6666

$ws.Range("H65").Value = 1975.4546
$ws.Range("I65").Value = 1966.375
$ws.Range("J65").Value = 1999.6666
$ws.Range("K65").Value = 9831.875
$ws.Range("L65").Value = 9998.333000000001
$ws.Range("M65").Value = -6711.875
$ws.Range("N65").Value = -16238.333

$ws.Range("H100").Value = 7620.0527
$ws.Range("I100").Value = 6860.875
$ws.Range("J100").Value = 8172.1816
$ws.Range("K100").Value = 6860.875
$ws.Range("L100").Value = 8172.1816
$ws.Range("M100").Value = -6319.875
$ws.Range("N100").Value = -9254.1816

$ws.Range("H106").Value = 7769.4614
$ws.Range("I106").Value = 7554.8184
$ws.Range("J106").Value = 8950
$ws.Range("K106").Value = 7554.8184
$ws.Range("L106").Value = 8950
$ws.Range("M106").Value = -6923.8184
$ws.Range("N106").Value = -10212

$ws.Range("H113").Value = 9518.904
$ws.Range("I113").Value = 8444.111000000001
$ws.Range("J113").Value = 10325
$ws.Range("K113").Value = 8444.111000000001
$ws.Range("L113").Value = 10325
$ws.Range("M113").Value = -5190.111000000001
$ws.Range("N113").Value = -16833

$ws.Range("H138").Value = 2903.5967
$ws.Range("I138").Value = 1700.9565
$ws.Range("J138").Value = 3612.8462
$ws.Range("K138").Value = 5102.8695
$ws.Range("L138").Value = 10838.5386
$ws.Range("M138").Value = 37.13050000000021
$ws.Range("N138").Value = -21118.5386

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 7515.362
$ws.Range("I32").Value = 8173.7617
$ws.Range("J32").Value = 1984.8
$ws.Range("K32").Value = 8173.7617
$ws.Range("L32").Value = 1984.8
$ws.Range("M32").Value = -7886.7617
$ws.Range("N32").Value = -2558.8

$ws.Range("H35").Value = 18000
$ws.Range("I35").Value = 1000
$ws.Range("J35").Value = 35000
$ws.Range("K35").Value = 1000
$ws.Range("L35").Value = 35000
$ws.Range("M35").Value = -594
$ws.Range("N35").Value = -35812

$ws.Range("H45").Value = 3366
$ws.Range("I45").Value = 2357.25
$ws.Range("J45").Value = 4374.75
$ws.Range("K45").Value = 2357.25
$ws.Range("L45").Value = 4374.75
$ws.Range("M45").Value = -1980.25
$ws.Range("N45").Value = -5128.75

$ws.Range("H60").Value = 16500
$ws.Range("I60").Value = 8000
$ws.Range("J60").Value = 25000
$ws.Range("K60").Value = 8000
$ws.Range("L60").Value = 25000
$ws.Range("M60").Value = -7267
$ws.Range("N60").Value = -26466

$ws.Range("H97").Value = 14499.667
$ws.Range("I97").Value = 19999
$ws.Range("J97").Value = 11750
$ws.Range("K97").Value = 19999
$ws.Range("L97").Value = 11750
$ws.Range("M97").Value = -19503
$ws.Range("N97").Value = -12742

$ws.Range("H98").Value = 55750
$ws.Range("I98").Value = 0
$ws.Range("J98").Value = 55750
$ws.Range("K98").Value = 0
$ws.Range("L98").Value = 55750
$ws.Range("N98").Value = -61740

$ws.Range("H102").Value = 2947.353
$ws.Range("I102").Value = 2947.353
$ws.Range("J102").Value = 0
$ws.Range("K102").Value = 2947.353
$ws.Range("L102").Value = 0
$ws.Range("M102").Value = -1325.353

$ws.Range("H111").Value = 45000
$ws.Range("I111").Value = 0
$ws.Range("J111").Value = 45000
$ws.Range("K111").Value = 0
$ws.Range("L111").Value = 45000
$ws.Range("M111").ClearContents()
$ws.Range("N111").Value = -53180

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value = 12501760
$ws.Range("I86").Value = 14707190
$ws.Range("J86").Value = 4321.3335
$ws.Range("K86").Value = 14707190
$ws.Range("L86").Value = 4321.3335
$ws.Range("M86").Value = -14706067
$ws.Range("N86").Value = -6567.3335

$ws.Range("H89").Value = 12501760
$ws.Range("I89").Value = 14707190
$ws.Range("J89").Value = 4321.3335
$ws.Range("K89").Value = 73535950
$ws.Range("L89").Value = 21606.6675
$ws.Range("M89").Value = -73530334
$ws.Range("N89").Value = -32838.6675

$ws.Range("H99").Value = 200905
$ws.Range("I99").Value = 200905
$ws.Range("J99").Value = 0
$ws.Range("K99").Value = 200905
$ws.Range("L99").Value = 0
$ws.Range("M99").Value = -199407

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H7").Value = 237.86667
$ws.Range("I7").Value = 120.85714
$ws.Range("J7").Value = 340.25
$ws.Range("K7").Value = 120.85714
$ws.Range("L7").Value = 340.25
$ws.Range("M7").Value = -7.857140000000001
$ws.Range("N7").Value = -566.25

$ws.Range("H59").Value = 514999.5
$ws.Range("I59").Value = 514999.5
$ws.Range("J59").Value = 0
$ws.Range("K59").Value = 514999.5
$ws.Range("L59").Value = 0
$ws.Range("M59").Value = -513854.5

$ws.Range("H60").Value = 398.26666
$ws.Range("I60").Value = 398.26666
$ws.Range("J60").Value = 0
$ws.Range("K60").Value = 398.26666
$ws.Range("L60").Value = 0
$ws.Range("M60").Value = 112.73334

$ws.Range("H88").Value = 22781
$ws.Range("I88").Value = 0
$ws.Range("J88").Value = 22781
$ws.Range("K88").Value = 0
$ws.Range("L88").Value = 22781
$ws.Range("N88").Value = -23593

$ws.Range("H91").Value = 22781
$ws.Range("I91").Value = 0
$ws.Range("J91").Value = 22781
$ws.Range("K91").Value = 0
$ws.Range("L91").Value = 22781
$ws.Range("N91").Value = -25589

$ws.Range("H105").Value = 839
$ws.Range("I105").Value = 790.9
$ws.Range("J105").Value = 999.3333
$ws.Range("K105").Value = 790.9
$ws.Range("L105").Value = 999.3333
$ws.Range("M105").Value = 956.1
$ws.Range("N105").Value = -4493.3333

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H26").Value = 564.7778
$ws.Range("I26").Value = 14.6
$ws.Range("J26").Value = 1252.5
$ws.Range("K26").Value = 43.8
$ws.Range("L26").Value = 3757.5
$ws.Range("M26").Value = 244.2
$ws.Range("N26").Value = -4333.5

$ws.Range("H109").Value = 398.5
$ws.Range("I109").Value = 398.5
$ws.Range("J109").Value = 0
$ws.Range("K109").Value = 1195.5
$ws.Range("L109").Value = 0
$ws.Range("M109").Value = -155.5

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H46").Value = 1000
$ws.Range("I46").Value = 1000
$ws.Range("J46").Value = 0
$ws.Range("K46").Value = 1000
$ws.Range("L46").Value = 0
$ws.Range("M46").Value = -844

$ws.Range("H53").Value = 14000
$ws.Range("I53").Value = 0
$ws.Range("J53").Value = 14000
$ws.Range("K53").Value = 0
$ws.Range("L53").Value = 14000
$ws.Range("N53").Value = -15262

$ws.Range("H100").Value = 39160.707
$ws.Range("I100").Value = 0
$ws.Range("J100").Value = 39160.707
$ws.Range("K100").Value = 0
$ws.Range("L100").Value = 39160.707
$ws.Range("N100").Value = -41324.707

$ws.Range("H122").Value = 9224.625
$ws.Range("I122").Value = 7666.6665
$ws.Range("J122").Value = 10159.4
$ws.Range("K122").Value = 22999.9995
$ws.Range("L122").Value = 30478.2
$ws.Range("M122").Value = -20549.9995
$ws.Range("N122").Value = -35378.2

$ws.Range("H132").Value = 6436.8286
$ws.Range("I132").Value = 5805.9062
$ws.Range("J132").Value = 13166.667
$ws.Range("K132").Value = 17417.7186
$ws.Range("L132").Value = 39500.001
$ws.Range("M132").Value = -14887.7186
$ws.Range("N132").Value = -44560.001

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H46").Value = 2230.1428
$ws.Range("I46").Value = 1231.6666
$ws.Range("J46").Value = 2629.5334
$ws.Range("K46").Value = 1231.6666
$ws.Range("L46").Value = 2629.5334
$ws.Range("M46").Value = -1043.6666
$ws.Range("N46").Value = -3005.5334

$ws.Range("H55").Value = 749.8333
$ws.Range("I55").Value = 780
$ws.Range("J55").Value = 689.5
$ws.Range("K55").Value = 780
$ws.Range("L55").Value = 689.5
$ws.Range("M55").Value = -607
$ws.Range("N55").Value = -1035.5

$ws.Range("H100").Value = 58756
$ws.Range("I100").Value = 105488.73
$ws.Range("J100").Value = 7350
$ws.Range("K100").Value = 105488.73
$ws.Range("L100").Value = 7350
$ws.Range("M100").Value = -104947.73
$ws.Range("N100").Value = -8432

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H96").Value = 53636
$ws.Range("I96").Value = 103371.8
$ws.Range("J96").Value = 3900.2
$ws.Range("K96").Value = 103371.8
$ws.Range("L96").Value = 3900.2
$ws.Range("M96").Value = -101998.8
$ws.Range("N96").Value = -6646.2

$ws.Range("H122").Value = 2178.8
$ws.Range("I122").Value = 1998
$ws.Range("J122").Value = 2450
$ws.Range("K122").Value = 5994
$ws.Range("L122").Value = 7350
$ws.Range("M122").Value = -3544
$ws.Range("N122").Value = -12250

